$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary row fixes (No. / Marking / Total), rows 10-12
# ---------------------------------------------------------------------------

# Give A10, A11, A12 the same "header" formatting already used by A9
# (bold/centered/bordered), without touching their existing text values.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

# Row 10 ("No."): Right / Not Attempt / Max
$ws.Range("B10").Value = 14
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 28

# Row 11 ("Marking"): Right / Wrong (now a real negative number, not text)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total"): Right score / fraction text
$ws.Range("B12").Value = 56
$ws.Range("E12").Value = "56/112"

# ---------------------------------------------------------------------------
# Answer-key table: only one Student/Correct Ans pair (columns A & B) is
# kept now, so the duplicate pairs in D:E and G:H are dropped.
# ---------------------------------------------------------------------------

# Drop the 3rd pair's header ("Student Ans" / "Correct Ans" in G15:H15)
# together with every data cell below it (G16:H21), and drop the 2nd
# pair's data for every question row beyond the first three (D19:E40).
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# Fill in the "Student Ans" column (A) for every question the student
# answered correctly - it now mirrors the "Correct Ans" column (B), shown
# in the "right answer" (green) style already used by B10:B12.
$ws.Range("B10").Copy()
$correctRows = @(16, 17, 19, 22, 24, 25, 27, 29, 30, 32, 33, 35, 37, 39)
foreach ($r in $correctRows) {
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $ws.Range("B$r").Text
}
